$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "September 19, 2025*") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# ---------------------------------------------------------------------
# 2) Mailing address: split the single line
#       "1730 Highland Place, Berkeley CA 94709"
#    (the one in the letter header, NOT the copy inside the reply-slip
#    table) into two paragraphs:
#       "1730 Highland Place"
#       "Berkeley, CA 94709"
#    using InsertXML so the new run keeps xml:space="preserve" and the
#    exact run/paragraph formatting used elsewhere in the letter.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if (($p.Range.Text -like "1730 Highland Place, Berkeley CA 94709*") -and ($p.Range.Information(12) -eq $false)) {
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
          '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
          '<pkg:xmlData>' + `
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
          '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">1730 Highland Place</w:t></w:r></w:p>' + `
          '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Berkeley, CA 94709</w:t></w:r></w:p>' + `
          '</w:body></w:document>' + `
          '</pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------
# 3) Remove the now-superfluous blank "NoSpacing" paragraph that sits
#    right after the "... Board of Directors" line.
# ---------------------------------------------------------------------
$foundBoard = $false
$blankAfterBoard = $null
foreach ($p in $d.Paragraphs) {
    if ($foundBoard) {
        $blankAfterBoard = $p
        $foundBoard = $false
    }
    if ($p.Range.Text -like "*Board of Directors*") {
        $foundBoard = $true
    }
}
if ($blankAfterBoard -ne $null) {
    $blankAfterBoard.Range.Delete()
}
